$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1780821917808219
$ws.Range("C2").Value = 0.54337899543379
$ws.Range("J2").Value = 0.0182648401826484
$ws.Range("P2").Value = 0.1278538812785388
$ws.Range("S2").Value = 0.1324200913242009
$ws.Range("B3").Value = 0.008474576271186441
$ws.Range("J3").Value = 0.03389830508474576
$ws.Range("P3").Value = 0.7288135593220338
$ws.Range("S3").Value = 0.2288135593220339
$ws.Range("J4").Value = 0.08333333333333333
$ws.Range("P4").Value = 0.7222222222222222
$ws.Range("S4").Value = 0.1944444444444444
$ws.Range("S5").Value = 1
$ws.Range("B6").Value = 0.05911330049261083
$ws.Range("D6").Value = 0.009852216748768473
$ws.Range("F6").Value = 0.06896551724137931
$ws.Range("J6").Value = 0.2315270935960591
$ws.Range("O6").Value = 0.009852216748768473
$ws.Range("Q6").Value = 0.1280788177339902
$ws.Range("R6").Value = 0.05911330049261083
$ws.Range("S6").Value = 0.4334975369458128
$ws.Range("B7").Value = 0.07766990291262135
$ws.Range("D7").Value = 0.01456310679611651
$ws.Range("E7").Value = 0.004854368932038835
$ws.Range("F7").Value = 0.04368932038834952
$ws.Range("J7").Value = 0.1504854368932039
$ws.Range("O7").Value = 0.009708737864077669
$ws.Range("Q7").Value = 0.1310679611650485
$ws.Range("R7").Value = 0.06310679611650485
$ws.Range("S7").Value = 0.5048543689320388
$ws.Range("B8").Value = 0.118895966029724
$ws.Range("D8").Value = 0.01698513800424628
$ws.Range("F8").Value = 0.06581740976645435
$ws.Range("J8").Value = 0.1380042462845011
$ws.Range("O8").Value = 0.0148619957537155
$ws.Range("Q8").Value = 0.1528662420382166
$ws.Range("R8").Value = 0.06794055201698514
$ws.Range("S8").Value = 0.4246284501061571
$ws.Range("B9").Value = 0.119047619047619
$ws.Range("D9").Value = 0.01785714285714286
$ws.Range("F9").Value = 0.09523809523809523
$ws.Range("J9").Value = 0.119047619047619
$ws.Range("Q9").Value = 0.09523809523809523
$ws.Range("R9").Value = 0.09523809523809523
$ws.Range("S9").Value = 0.4583333333333333
$ws.Range("B10").Value = 0.09455842997323818
$ws.Range("D10").Value = 0.01784121320249777
$ws.Range("F10").Value = 0.05798394290811775
$ws.Range("J10").Value = 0.1275646743978591
$ws.Range("O10").Value = 0.008028545941123996
$ws.Range("Q10").Value = 0.223907225691347
$ws.Range("R10").Value = 0.09099018733273863
$ws.Range("S10").Value = 0.3791257805530776
$ws.Range("G11").Value = 0.1851851851851852
$ws.Range("J11").Value = 0.07777777777777778
$ws.Range("K11").Value = 0.1962962962962963
$ws.Range("L11").Value = 0.5148148148148148
$ws.Range("S11").Value = 0.02592592592592593
$ws.Range("G12").Value = 0.7554347826086957
$ws.Range("J12").Value = 0.1576086956521739
$ws.Range("K12").Value = 0.0108695652173913
$ws.Range("L12").Value = 0.03260869565217391
$ws.Range("S12").Value = 0.04347826086956522
$ws.Range("G13").Value = 0.5714285714285714
$ws.Range("J13").Value = 0.3095238095238095
$ws.Range("S13").Value = 0.119047619047619
$ws.Range("F15").Value = 0.01470588235294118
$ws.Range("H15").Value = 0.1691176470588235
$ws.Range("I15").Value = 0.08823529411764706
$ws.Range("J15").Value = 0.3602941176470588
$ws.Range("K15").Value = 0.05882352941176471
$ws.Range("M15").Value = 0.01470588235294118
$ws.Range("O15").Value = 0.03676470588235294
$ws.Range("S15").Value = 0.2573529411764706
$ws.Range("F16").Value = 0.01438848920863309
$ws.Range("H16").Value = 0.1870503597122302
$ws.Range("I16").Value = 0.06474820143884892
$ws.Range("J16").Value = 0.4028776978417266
$ws.Range("K16").Value = 0.1438848920863309
$ws.Range("M16").Value = 0.007194244604316547
$ws.Range("N16").Value = 0.007194244604316547
$ws.Range("O16").Value = 0.02877697841726619
$ws.Range("S16").Value = 0.1438848920863309
$ws.Range("F17").Value = 0.02083333333333333
$ws.Range("H17").Value = 0.2135416666666667
$ws.Range("I17").Value = 0.09635416666666667
$ws.Range("J17").Value = 0.4479166666666667
$ws.Range("K17").Value = 0.09375
$ws.Range("M17").Value = 0.01822916666666667
$ws.Range("O17").Value = 0.04166666666666666
$ws.Range("S17").Value = 0.06770833333333333
$ws.Range("F18").Value = 0.02298850574712644
$ws.Range("H18").Value = 0.1954022988505747
$ws.Range("I18").Value = 0.07471264367816093
$ws.Range("J18").Value = 0.4310344827586207
$ws.Range("K18").Value = 0.1379310344827586
$ws.Range("M18").Value = 0.02298850574712644
$ws.Range("N18").Value = 0.01149425287356322
$ws.Range("O18").Value = 0.04022988505747126
$ws.Range("S18").Value = 0.06321839080459771
$ws.Range("F19").Value = 0.02931596091205212
$ws.Range("H19").Value = 0.2255700325732899
$ws.Range("I19").Value = 0.0741042345276873
$ws.Range("J19").Value = 0.3509771986970684
$ws.Range("K19").Value = 0.08794788273615635
$ws.Range("M19").Value = 0.02117263843648208
$ws.Range("O19").Value = 0.0504885993485342
$ws.Range("S19").Value = 0.1604234527687296
